$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove asterisk from "Business Criticality*" header label in row 2
$ws.Range("B2").Value = "Business Criticality"

# Add a new "Mandatory*" label above column B, styled bold + wrap text
$ws.Range("B1").Value = "Mandatory*"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").WrapText = $true

# Remove asterisk from "Application Name*" header label in row 2
$ws.Range("A2").Value = "Application Name"

# Update selection
$ws.Range("B6").Select()
